$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2149751.8
$ws.Range("J17").Value = 2315079
$ws.Range("L17").Value = 6945237
$ws.Range("N17").Value = -6945573
$ws.Range("H92").Value = 1378.8667
$ws.Range("I92").Value = 1857.6364
$ws.Range("K92").Value = 1857.6364
$ws.Range("M92").Value = -609.6364000000001
$ws.Range("H98").Value = 2869.8333
$ws.Range("I98").Value = 3034.8125
$ws.Range("K98").Value = 3034.8125
$ws.Range("M98").Value = -1536.8125
$ws.Range("H122").Value = 2869.8333
$ws.Range("I122").Value = 3034.8125
$ws.Range("K122").Value = 9104.4375
$ws.Range("M122").Value = -6654.4375
$ws.Range("H125").Value = 6042.1113
$ws.Range("I125").Value = 6109.8335
$ws.Range("J125").Value = 5906.6665
$ws.Range("K125").Value = 54988.5015
$ws.Range("L125").Value = 53159.9985
$ws.Range("M125").Value = -52528.5015
$ws.Range("N125").Value = -58079.9985
$ws.Range("H131").Value = 2342.8
$ws.Range("I131").Value = 679.125
$ws.Range("K131").Value = 2037.375
$ws.Range("M131").Value = 3002.625
$ws.Range("H132").Value = 12540.186
$ws.Range("I132").Value = 11334.315
$ws.Range("K132").Value = 34002.945
$ws.Range("M132").Value = -31472.945
$ws.Range("H135").Value = 3907.875
$ws.Range("I135").Value = 2346.2666
$ws.Range("K135").Value = 21116.3994
$ws.Range("M135").Value = -18581.3994
$ws.Range("H137").Value = 7686.9453
$ws.Range("J137").Value = 16659.592
$ws.Range("L137").Value = 49978.776
$ws.Range("N137").Value = -55078.776
$ws.Range("H138").Value = 3765.8262
$ws.Range("I138").Value = 6367.875
$ws.Range("J138").Value = 2378.0667
$ws.Range("K138").Value = 19103.625
$ws.Range("L138").Value = 7134.2001
$ws.Range("M138").Value = -13963.625
$ws.Range("N138").Value = -17414.2001
$ws.Range("H141").Value = 4005.7144
$ws.Range("I141").Value = 3658
$ws.Range("J141").Value = 4875
$ws.Range("K141").Value = 10974
$ws.Range("L141").Value = 14625
$ws.Range("M141").Value = -5794
$ws.Range("N141").Value = -24985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3263.1035
$ws.Range("I2").Value = 1313.7778
$ws.Range("K2").Value = 1313.7778
$ws.Range("M2").Value = -1200.7778
$ws.Range("H32").Value = 6240.07
$ws.Range("I32").Value = 4584.35
$ws.Range("K32").Value = 4584.35
$ws.Range("M32").Value = -4297.35
$ws.Range("H45").Value = 4057.1428
$ws.Range("I45").Value = 5550
$ws.Range("J45").Value = 3460
$ws.Range("K45").Value = 5550
$ws.Range("L45").Value = 3460
$ws.Range("M45").Value = -5173
$ws.Range("N45").Value = -4214
$ws.Range("H110").Value = 3245.3845
$ws.Range("I110").Value = 2148.625
$ws.Range("K110").Value = 2148.625
$ws.Range("M110").Value = -103.625
$ws.Range("H116").Value = 3263.1035
$ws.Range("I116").Value = 1313.7778
$ws.Range("K116").Value = 1313.7778
$ws.Range("M116").Value = 980.2221999999999
$ws.Range("H122").Value = 2625.7144
$ws.Range("I122").Value = 2507
$ws.Range("K122").Value = 7521
$ws.Range("M122").Value = -5071
$ws.Range("H126").Value = 8332.666999999999
$ws.Range("I126").Value = 8332.666999999999
$ws.Range("K126").Value = 24998.001
$ws.Range("M126").Value = -22528.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3263.1035
$ws.Range("I3").Value = 1313.7778
$ws.Range("K3").Value = 1313.7778
$ws.Range("M3").Value = -1199.7778
$ws.Range("H20").Value = 14716.866
$ws.Range("I20").Value = 4718.485
$ws.Range("J20").Value = 26937.111
$ws.Range("K20").Value = 4718.485
$ws.Range("L20").Value = 26937.111
$ws.Range("M20").Value = -4471.485
$ws.Range("N20").Value = -27431.111
$ws.Range("H33").Value = 7000
$ws.Range("J33").Value = 7000
$ws.Range("L33").Value = 7000
$ws.Range("N33").Value = -7672
$ws.Range("H94").Value = 3486.3044
$ws.Range("I94").Value = 3601.6667
$ws.Range("K94").Value = 3601.6667
$ws.Range("M94").Value = -3150.6667
$ws.Range("H105").Value = 1672.9333
$ws.Range("I105").Value = 1093.5555
$ws.Range("K105").Value = 1093.5555
$ws.Range("M105").Value = 653.4445000000001
$ws.Range("H128").Value = 7008.9
$ws.Range("I128").Value = 7008.9
$ws.Range("K128").Value = 21026.7
$ws.Range("M128").Value = -18536.7
$ws.Range("H134").Value = 6751.2144
$ws.Range("I134").Value = 1598.0286
$ws.Range("J134").Value = 32517.143
$ws.Range("K134").Value = 4794.085800000001
$ws.Range("L134").Value = 97551.429
$ws.Range("M134").Value = -2259.085800000001
$ws.Range("N134").Value = -102621.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24607.637
$ws.Range("I31").Value = 8180.1875
$ws.Range("J31").Value = 68414.164
$ws.Range("K31").Value = 8180.1875
$ws.Range("L31").Value = 68414.164
$ws.Range("M31").Value = -7885.1875
$ws.Range("N31").Value = -69004.164
$ws.Range("H34").Value = 24607.637
$ws.Range("I34").Value = 8180.1875
$ws.Range("J34").Value = 68414.164
$ws.Range("K34").Value = 8180.1875
$ws.Range("L34").Value = 68414.164
$ws.Range("M34").Value = -7978.1875
$ws.Range("N34").Value = -68818.164
$ws.Range("H134").Value = 29418322
$ws.Range("I134").Value = 1482.2106
$ws.Range("J134").Value = 66679652
$ws.Range("K134").Value = 4446.6318
$ws.Range("L134").Value = 200038956
$ws.Range("M134").Value = -1911.6318
$ws.Range("N134").Value = -200044026

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47091660
$ws.Range("I4").Value = 51700824
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 155102472
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -155102360
$ws.Range("N4").Value = -3000224
$ws.Range("H98").Value = 12058
$ws.Range("J98").Value = 12058
$ws.Range("L98").Value = 36174
$ws.Range("N98").Value = -39170
$ws.Range("H129").Value = 3000
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000
$ws.Range("H131").Value = 1499.54
$ws.Range("J131").Value = 1499.54
$ws.Range("L131").Value = 4498.62
$ws.Range("N131").Value = -14578.62
$ws.Range("H137").Value = 1457.6666
$ws.Range("J137").Value = 1956.1428
$ws.Range("L137").Value = 5868.428400000001
$ws.Range("N137").Value = -16068.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21406.5
$ws.Range("I70").Value = 4994
$ws.Range("K70").Value = 4994
$ws.Range("M70").Value = -4724
$ws.Range("H73").Value = 21406.5
$ws.Range("I73").Value = 4994
$ws.Range("K73").Value = 4994
$ws.Range("M73").Value = -4058
$ws.Range("H122").Value = 3258.8333
$ws.Range("I122").Value = 3155.0908
$ws.Range("J122").Value = 4400
$ws.Range("K122").Value = 9465.2724
$ws.Range("L122").Value = 13200
$ws.Range("M122").Value = -7015.2724
$ws.Range("N122").Value = -18100
$ws.Range("H132").Value = 30157.5
$ws.Range("I132").Value = 22541
$ws.Range("J132").Value = 53007
$ws.Range("K132").Value = 67623
$ws.Range("L132").Value = 159021
$ws.Range("M132").Value = -65093
$ws.Range("N132").Value = -164081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 4500450
$ws.Range("I23").Value = 9000000
$ws.Range("J23").Value = 900
$ws.Range("K23").Value = 9000000
$ws.Range("L23").Value = 900
$ws.Range("M23").Value = -8999770
$ws.Range("N23").Value = -1360
$ws.Range("H55").Value = 2033.4878
$ws.Range("I55").Value = 982.0454999999999
$ws.Range("K55").Value = 982.0454999999999
$ws.Range("M55").Value = -809.0454999999999
$ws.Range("H61").Value = 3210.5881
$ws.Range("I61").Value = 2089.818
$ws.Range("K61").Value = 2089.818
$ws.Range("M61").Value = -1887.818
$ws.Range("H82").Value = 4689.3
$ws.Range("I82").Value = 3968.2104
$ws.Range("J82").Value = 5934.8184
$ws.Range("K82").Value = 3968.2104
$ws.Range("L82").Value = 5934.8184
$ws.Range("M82").Value = -3607.2104
$ws.Range("N82").Value = -6656.8184
$ws.Range("H85").Value = 4689.3
$ws.Range("I85").Value = 3968.2104
$ws.Range("J85").Value = 5934.8184
$ws.Range("K85").Value = 3968.2104
$ws.Range("L85").Value = 5934.8184
$ws.Range("M85").Value = -2720.2104
$ws.Range("N85").Value = -8430.8184
$ws.Range("H107").Value = 7439.706
$ws.Range("I107").Value = 7439.706
$ws.Range("K107").Value = 7439.706
$ws.Range("M107").Value = -5519.706
$ws.Range("H113").Value = 3210.5881
$ws.Range("I113").Value = 2089.818
$ws.Range("K113").Value = 2089.818
$ws.Range("M113").Value = 80.18199999999979
$ws.Range("H122").Value = 7794.25
$ws.Range("I122").Value = 5404.5713
$ws.Range("K122").Value = 16213.7139
$ws.Range("M122").Value = -13763.7139
$ws.Range("H132").Value = 1917408.6
$ws.Range("I132").Value = 3158.5
$ws.Range("J132").Value = 4469742
$ws.Range("K132").Value = 9475.5
$ws.Range("L132").Value = 13409226
$ws.Range("M132").Value = -6945.5
$ws.Range("N132").Value = -13414286
$ws.Range("H136").Value = 19242.2
$ws.Range("J136").Value = 15941.625
$ws.Range("L136").Value = 47824.875
$ws.Range("N136").Value = -52924.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 52000
$ws.Range("I47").Value = 52000
$ws.Range("K47").Value = 52000
$ws.Range("M47").Value = -51428
$ws.Range("H122").Value = 4874.909
$ws.Range("I122").Value = 2115.5715
$ws.Range("J122").Value = 9703.75
$ws.Range("K122").Value = 6346.7145
$ws.Range("L122").Value = 29111.25
$ws.Range("M122").Value = -3896.7145
$ws.Range("N122").Value = -34011.25

Write-Output "Applied 247 cell updates across 8 sheets"